$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-only cells (Price/Volume columns carry textual data, e.g.
# "29.421.61" or "0.9990", that Excel would otherwise auto-convert to
# numbers/dates). Setting NumberFormat to "@" before the assignment keeps
# the literal text intact, matching the original inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.421.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7181"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07825"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3100"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.83"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08257"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.874.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7270"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.284"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.380.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.912"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "245.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007902"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.904"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.98"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.024"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.487"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.397"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.143"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05286"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7226"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.676"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01869"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.234.65"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.719"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9098"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.42"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.088"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5332"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.933"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +12.09%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.757"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4333"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.257"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.079"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.28%  "
